# Add the new bathtub-product test-case row beneath the existing SKU header
# row, then move the active selection down to the next empty row (A3) —
# matching where Excel leaves the cursor after data entry on row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 105546
$ws.Range("B2").Value = 139398

[void]$ws.Range("A3").Select()
